$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update revised figures for existing rows 129-132 ---

# Row 129
$ws.Range("C129").Value = 298.6
$ws.Range("D129").Value = 1236.6
$ws.Range("G129").Value = 1521.6
$ws.Range("I129").Value = 964.9
$ws.Range("K129").Value = 1352.6

# Row 130
$ws.Range("C130").Value = 312
$ws.Range("D130").Value = 1278
$ws.Range("E130").Value = 994.5
$ws.Range("F130").Value = 468.9
$ws.Range("G130").Value = 1553.9
$ws.Range("H130").Value = 201.7
$ws.Range("I130").Value = 1001.4
$ws.Range("J130").Value = 636.5
$ws.Range("K130").Value = 1426.4

# Row 131
$ws.Range("C131").Value = 299.9
$ws.Range("D131").Value = 1252.6
$ws.Range("E131").Value = 997.8
$ws.Range("F131").Value = 462.5
$ws.Range("G131").Value = 1626.4
$ws.Range("I131").Value = 1007.9
$ws.Range("J131").Value = 634.1

# Row 132
$ws.Range("C132").Value = 305.8
$ws.Range("D132").Value = 1253.5
$ws.Range("E132").Value = 1013.8
$ws.Range("F132").Value = 513.5
$ws.Range("H132").Value = 199.7
$ws.Range("I132").Value = 999.8
$ws.Range("J132").Value = 617.1
$ws.Range("K132").Value = 1529.6

# --- Append new row 140 with the latest monthly data (01-09-2021) ---

$ws.Range("A140").Value = "'01-09-2021"
$ws.Range("A140").Style = "Normal"
$ws.Range("B140").Value = 8345.200000000001
$ws.Range("C140").Value = 335.7
$ws.Range("D140").Value = 1281.6
$ws.Range("E140").Value = 1035.5
$ws.Range("F140").Value = 478.7
$ws.Range("G140").Value = 1740.4
$ws.Range("H140").Value = 187
$ws.Range("I140").Value = 1100.8
$ws.Range("J140").Value = 649.9
$ws.Range("K140").Value = 1493
$ws.Range("L140").Value = 26.5
$ws.Range("M140").Value = 16.1
